# Mise à jour de l'application : ajout de la séance d'entrainement du 11/10/2025
# (nouvelle colonne BJ) + saisie des présences de chaque joueur pour cette séance,
# et correction de 3 présences mal saisies pour Jassim Assoul.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) New training date header (11 Oct 2025 -> serial 45941)
$ws.Range("BJ1").Value = 45941

# 2) Attendance codes entered for the new session, per player row
#    (row 12 - Yanis Berrached - has no entry for this session)
$attendance = @{
    2  = "P"
    3  = "R"
    4  = "P"
    5  = "P"
    6  = "REP"
    7  = "P"
    8  = "B"
    9  = "P"
    10 = "P"
    11 = "P"
    13 = "B"
    14 = "P"
    15 = "P"
    16 = "P"
    17 = "B"
    18 = "P"
    19 = "P"
    20 = "P"
    21 = "M"
    22 = "P"
    23 = "B"
    24 = "P"
    25 = "P"
    26 = "P"
    27 = "P"
    28 = "P"
    29 = "P"
}

foreach ($row in $attendance.Keys) {
    $ws.Range("BJ$row").Value = $attendance[$row]
}

# 3) Correction of three previously mis-entered codes for Jassim Assoul (row 3):
#    the BG/BH/BI sessions were actually "P", not "R".
$ws.Range("BG3").Value = "P"
$ws.Range("BH3").Value = "P"
$ws.Range("BI3").Value = "P"

# 4) Extend formatting from the last existing date column (BI) to the new one (BJ),
#    skipping row 12 which has no data in column BI. Values were already written
#    above, so this paste (format-only) won't disturb them.
$ws.Range("BI1:BI11").Copy() | Out-Null
$ws.Range("BJ1:BJ11").PasteSpecial(-4122) | Out-Null

$ws.Range("BI13:BI29").Copy() | Out-Null
$ws.Range("BJ13:BJ29").PasteSpecial(-4122) | Out-Null

$ws.Application.CutCopyMode = $false

# 5) Recalculate all formulas (attendance totals per player)
$excel.CalculateFull()

# 6) Restore the view: frozen pane boundary + current selection
$ws.Range("BD1").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("BM4").Select() | Out-Null
